# Fruta / hortaliza, semanal
# Insert a new weekly record at row 245 of the data table (pushing existing
# rows 245:356 down to 246:357) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 245, shifting rows 245-356 down to 246-357.
$ws.Rows("245:245").Insert()

# Populate the newly inserted row 245 with the new data point.
$ws.Cells.Item(245, 1).Value2  = 10
$ws.Cells.Item(245, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(245, 3).Value2  = "La Araucanía"
$ws.Cells.Item(245, 4).Value2  = 44839
$ws.Cells.Item(245, 5).Value2  = 9
$ws.Cells.Item(245, 6).Value2  = 100112001
$ws.Cells.Item(245, 7).Value2  = "Berenjena"
$ws.Cells.Item(245, 8).Value2  = "Sin especificar"
$ws.Cells.Item(245, 9).Value2  = "Primera"
$ws.Cells.Item(245, 10).Value2 = 40
$ws.Cells.Item(245, 11).Value2 = 16000
$ws.Cells.Item(245, 12).Value2 = 16000
$ws.Cells.Item(245, 13).Value2 = 16000
$ws.Cells.Item(245, 14).Value2 = "`$/caja 40 unidades"
$ws.Cells.Item(245, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(245, 16).Value2 = 400
$ws.Cells.Item(245, 17).Value2 = 40
$ws.Cells.Item(245, 18).Value2 = "Hortaliza"
